# D2/Slides.pptx — "Added who will produce each graph/diagram"
#
# 1) Slide 5 "Class diagram so far": empty content placeholder gets
#    "<ed add class diagram>" (split "<" / "ed" / " add class diagram>").
# 2) Slide 6 "Architectural concepts": "<Diagram here>" becomes
#    "<Diagram here, ed>" (split "<Diagram " / "here, " / "ed" / ">").
# 3) Slide 8 "Time spent": the graph bullet gets " <kelvin>" tacked on
#    after the closing paren.
# 4) Slide 9 "Time spent": "Another graph" gets " <kelvin>" tacked on.

$p = $ppt.ActivePresentation

# --- Slide 5: Class diagram so far ------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange

$para5 = $tr5.Paragraphs(1, 1)
$para5.Text = "<ed add class diagram>"
# Touch the whole range so PowerPoint (re)attaches run properties to the
# freshly-typed text instead of leaving a bare, unformatted run.
$tr5.Text = $tr5.Text

$tr5.Characters(1, 1).Text = "<"
$tr5.Characters(2, 2).Text = "ed"

# --- Slide 6: Architectural concepts -----------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$tr6 = $shp6.TextFrame.TextRange

$para6 = $tr6.Paragraphs(3, 1)
$para6.Characters(10, 5).Text = "here, ed>"
$para6.Characters(16, 2).Text = "ed"

# --- Slide 8: Time spent (graph) ---------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$tr8 = $shp8.TextFrame.TextRange

$para8 = $tr8.Paragraphs(1, 1)
$para8.Characters(75, 1).Text = ") <kelvin>"

# --- Slide 9: Time spent (another graph) -------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tr9 = $shp9.TextFrame.TextRange

$para9 = $tr9.Paragraphs(1, 1)
$para9.Characters(9, 5).Text = "graph <kelvin>"
